$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

# Update Version value from 0.4.0 to 0.7.0
$wsMeta.Range("B3").Value = "0.7.0"

# Remove the "Jurisdiction" / "Chile" row entirely (row 11), shifting rows below up.
$wsMeta.Rows(11).Delete()
